$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Constants (standard Excel enum values) ---
# xlVAlignCenter = -4108, xlHAlignLeft = -4131

# ------------------------------------------------------------------
# The "Setting" / "Setting_Brightness" / "Setting_Volume" block that
# used to live at rows 5-7 is being moved further down the sheet
# (below the Player block), a blank separator row + a new "Popups"
# row are inserted, and the "Alert" row (previously row 12) moves
# down to make room.
# ------------------------------------------------------------------

# 1) Clear out the old "Setting" block (rows 5-7) - contents AND
#    formatting, since that whole block relocated further down.
$ws.Range("A5:F7").Clear()
$ws.Rows.Item(6).AutoFit()

# 2) Clear the old "Alert" row content (row 12) - it is being moved
#    down to row 14, leaving row 12 blank (still vertically centered).
$ws.Range("D12").Clear()
$ws.Range("A12").ClearContents()
$ws.Range("A12").VerticalAlignment = -4108
$ws.Rows.Item(12).AutoFit()
$ws.Range("A12").VerticalAlignment = -4108

# 3) New row 13: "Popups" label.
$ws.Range("A13").Value = "Popups"
$ws.Range("A13").VerticalAlignment = -4108

# 4) New row 14: "Alert" (moved down from row 12) with its Keypad
#    Enter note in column D.
$ws.Range("A14").Value = "Alert"
$ws.Range("A14").VerticalAlignment = -4108

$ws.Range("D14").Value = "Keypad" + [char]10 + "LV_KEY_ENTER"
$ws.Range("D14").WrapText = $true

# 5) New row 15: "Setting" (moved down from row 5).
$ws.Range("A15").Value = "Setting"
$ws.Range("A15").VerticalAlignment = -4108

$ws.Range("C15").WrapText = $true
$ws.Range("C15").VerticalAlignment = -4108

# 6) New row 16: "Setting_Brightness" (moved down from row 6).
$ws.Range("A16").Value = "Setting_Brightness"
$ws.Range("A16").VerticalAlignment = -4108

$ws.Range("B16").Value = "Button (LV_EVENT_SHORT_CLICKED)" + [char]10 + "EVT_BUTON_SETTING_CLICKED"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").HorizontalAlignment = -4131
$ws.Range("B16").VerticalAlignment = -4108

$ws.Range("C16").Value = "Button (LV_EVENT_SHORT_CLICKED)" + [char]10 + "EVT_BUTTON_BACK_CLICKED"
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4108

$ws.Range("E16").Value = "Keypad" + [char]10 + "LV_KEY_UP"
$ws.Range("E16").WrapText = $true

$ws.Range("F16").Value = "Keypad" + [char]10 + "LV_KEY_DOWN"
$ws.Range("F16").WrapText = $true

# 7) New row 17: "Setting_Volume" (moved down from row 7).
$ws.Range("A17").Value = "Setting_Volume"
$ws.Range("A17").VerticalAlignment = -4108

# 8) Row heights for the wrapped rows settle at 30 (vs. the old 32)
#    once re-laid-out; match that here for the affected rows.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# 9) Sheet view: zoom back to 100%, move the active selection to the
#    relocated "Alert" Keypad-Enter cell.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("D14").Select()

"done"
